$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70; rows 70-87 shift down to 71-88.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new data record.
$ws.Cells.Item(70, 1).Value = 7
$ws.Cells.Item(70, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(70, 3).Value = "Ñuble"
$ws.Cells.Item(70, 4).Value = 45209
$ws.Cells.Item(70, 4).NumberFormat = $ws.Cells.Item(71, 4).NumberFormat
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 100112026
$ws.Cells.Item(70, 7).Value = "Haba"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 50
$ws.Cells.Item(70, 11).Value = 14000
$ws.Cells.Item(70, 12).Value = 14000
$ws.Cells.Item(70, 13).Value = 14000
$ws.Cells.Item(70, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(70, 16).Value = 560
$ws.Cells.Item(70, 17).Value = 25
$ws.Cells.Item(70, 18).Value = "Hortaliza"
